{"js": "// Apply the \"Added many more features\" edits to the document body.\n// Each entry is an exact (old, new) text pair that must be replaced\n// via Find/Replace-style search + insertText(\"Replace\").\nconst replacements = [\n  [\n    \"Play Celtic Charm Fire Blaze Quattro for Free - Review\",\n    \"Play Celtic Charm Fire Blaze Quattro for Free\",\n  ],\n  [\n    \"Impressive graphics and design\",\n    \"Impeccable graphics and design\",\n  ],\n  [\n    \"Multiple bonus rounds and free spins\",\n    \"Exciting bonus rounds and free spins feature\",\n  ],\n  [\n    \"Four independent reels with 27 ways to win\",\n    \"Four independent reels with up to 729 ways to win\",\n  ],\n  [\n    \"Medium volatility and a good RTP rate\",\n    \"Relaxing Celtic-inspired music and theme\",\n  ],\n  [\n    \"Inactive reels during free spins\",\n    \"Limited appearance of wild card symbol on specific reels\",\n  ],\n  [\n    \"Maximum bet of 500 \\u20ac might be limiting for high rollers\",\n    \"Only five free spins awarded during the free spins feature\",\n  ],\n  [\n    \"Discover the exciting gameplay mechanics and features of Celtic Charm Fire Blaze Quattro. Play for free with multiple bonus rounds and free spins.\",\n    \"Read our review to learn about the gameplay, graphics, bonus rounds, and jackpots of Celtic Charm Fire Blaze Quattro. Play now for free!\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits via Find/Replace on the\n# whole document story, mirroring Word's Ctrl+H \"Replace All\" behaviour.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Longer / more specific strings first so shorter ones below cannot\n# accidentally match inside them.\nReplace-AllText \"Discover the exciting gameplay mechanics and features of Celtic Charm Fire Blaze Quattro. Play for free with multiple bonus rounds and free spins.\" \"Read our review to learn about the gameplay, graphics, bonus rounds, and jackpots of Celtic Charm Fire Blaze Quattro. Play now for free!\"\n\nReplace-AllText \"Play Celtic Charm Fire Blaze Quattro for Free - Review\" \"Play Celtic Charm Fire Blaze Quattro for Free\"\n\nReplace-AllText \"Impressive graphics and design\" \"Impeccable graphics and design\"\nReplace-AllText \"Multiple bonus rounds and free spins\" \"Exciting bonus rounds and free spins feature\"\nReplace-AllText \"Four independent reels with 27 ways to win\" \"Four independent reels with up to 729 ways to win\"\nReplace-AllText \"Medium volatility and a good RTP rate\" \"Relaxing Celtic-inspired music and theme\"\nReplace-AllText \"Inactive reels during free spins\" \"Limited appearance of wild card symbol on specific reels\"\nReplace-AllText \"Maximum bet of 500 \u20ac might be limiting for high rollers\" \"Only five free spins awarded during the free spins feature\"\n\nWrite-Output \"done\"\n"}
